# Update countries & provincias Spain
# Applies:
#   1) Four country-name reorderings (the underlying statistics row kept its
#      position in the sheet while two countries' names - and any freshly
#      refreshed statistics - traded places / sort order).
#   2) Refreshed COVID-19 statistics (Casos totales, Nuevos casos, Casos
#      activos, Recuperados, Casos criticos, Muertes) for the countries whose
#      numbers changed between the 17:30 and 18:47 data pulls.
#   3) The "Datos actualizados" timestamp string, from 17:30 to 18:47.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1) Country name swaps (column A). Each pair of adjacent rows exchanges the
#    country name that is shown in it.
# ---------------------------------------------------------------------------
$ws.Range("A86").Value  = "Grecia"
$ws.Range("A87").Value  = "Camerun"

$ws.Range("A160").Value = "Republica de Chipre"
$ws.Range("A161").Value = "Togo"

$ws.Range("A207").Value = "Nueva Caledonia"
$ws.Range("A208").Value = "Santa Lucia"

$ws.Range("A215").Value = "Montserrat"
$ws.Range("A216").Value = "Islas Malvinas"

# ---------------------------------------------------------------------------
# 2) Refreshed statistics: row => @{Col = NewValue}
# ---------------------------------------------------------------------------
$updates = @{
    4   = @{ B = 7735938; C = 13192; D = 4952562; E = 2567251; G = 303; H = 216125 }
    23  = @{ B = 329138;  C = 1581;  D = 288954;  E = 31575;   G = 56;  H = 8609 }
    27  = @{ B = 280744;  C = 3718;  D = 216438;  E = 62488;   G = 21;  H = 1818 }
    29  = @{ B = 172806;  C = 1483;  D = 145403;  E = 17865;   G = 8;   H = 9538 }
    38  = @{ B = 116148;  C = 777;   D = 92157;   E = 21832;   G = 10;  H = 2159 }
    48  = @{ B = 92643;   C = 2621;  D = 50600;   E = 41220;   G = 29;  H = 823 }
    86  = @{ B = 20947;   C = 406;   D = 9989;    E = 10534;   G = 4;   H = 424 }
    87  = @{ B = 20924;   D = 19764; E = 740 }
    100 = @{ B = 12794;   C = 210;   D = 8907;    E = 3697;    G = 2;   H = 190 }
    102 = @{ D = 8500;    E = 2203 }
    111 = @{ B = 9119;    C = 140;   D = 7900;    E = 1091 }
    141 = @{ B = 3715;    C = 56;    E = 835 }
    160 = @{ B = 1897;    C = 21;    D = 1369;    E = 504;     G = 1;   H = 24 }
    161 = @{ B = 1881;    D = 1410;  E = 422;     H = 49 }
    165 = @{ B = 1355;    C = 1;     D = 1245;    E = 28 }
    215 = @{ D = 12;      H = 1 }
    216 = @{ D = 13;      H = 0 }
}

foreach ($row in $updates.Keys) {
    $cols = $updates[$row]
    foreach ($col in $cols.Keys) {
        $ws.Range("$col$row").Value = $cols[$col]
    }
}

# ---------------------------------------------------------------------------
# 3) Timestamp update in A1.
# ---------------------------------------------------------------------------
$ws.Range("A1").Value = "Datos actualizados a 7 de Octubre de 2020 a las 18:47"
